$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: Behavior Tree blog post update
$ws.Range("D28").Value = "Behavior Tree :: Design principles"
$ws.Range("E28").Value = "https://ropiens.tistory.com/213"

# Row 32: Spark UDF -> Airflow SubDag blog post update
$ws.Range("D32").Value = "[Airflow] SubDag 개념, 장단점, 샘플 코드 (feat. ChatGPT)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/409"

# Row 51: history command -> asksim.co.kr blog post update
$ws.Range("D51").Value = "여러분과 소통하기 위해 질의응답 사이트(asksim.co.kr)를 개발했습니다!"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EC%97%AC%EB%9F%AC%EB%B6%84%EA%B3%BC-%EC%86%8C%ED%86%B5%ED%95%98%EA%B8%B0-%EC%9C%84%ED%95%B4-%EC%A7%88%EC%9D%98%EC%9D%91%EB%8B%B5-%EC%82%AC%EC%9D%B4%ED%8A%B8asksimcokr%EB%A5%BC-%EA%B0%9C%EB%B0%9C%ED%96%88%EC%8A%B5%EB%8B%88%EB%8B%A4"
